# Histogram_delimiters.xlsx update
# Adds discrete-channel (Wires/Grids Start/Stop) table to replace the old
# Wire/Grid Separation list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet so stale shared-strings / cells don't linger.
$ws.Cells.Clear()

# --- Column width -------------------------------------------------------
$ws.Columns("A").ColumnWidth = 11.25

# --- Header labels -------------------------------------------------------
$ws.Range("A1").Value = "Wires"
$ws.Range("C1").Value = "Grids"

$ws.Range("A2").Value = "Start"
$ws.Range("B2").Value = "Stop"
$ws.Range("C2").Value = "Start"
$ws.Range("D2").Value = "Stop"

# --- Data rows -------------------------------------------------------
$ws.Range("A3").Value = 445
$ws.Range("B3").Value = 775
$ws.Range("C3").Value = 870
$ws.Range("D3").Value = 1125
$ws.Range("G3").Value = 525
$ws.Range("H3").Value = 780

$ws.Range("A4").Value = 780
$ws.Range("B4").Value = 1120

$ws.Range("A5").Value = 1120
$ws.Range("B5").Value = 1457

$ws.Range("A6").Value = 1457
$ws.Range("B6").Value = 1805

$ws.Range("A7").Value = 1810
$ws.Range("B7").Value = 2150

# --- Alignment (applied in single range calls so styles de-dupe) --------
$ws.Range("A1:D1").HorizontalAlignment = -4108
$ws.Range("A2:D2").HorizontalAlignment = -4108
$ws.Range("A2:D2").VerticalAlignment = -4108

# --- Merge the two group headers -----------------------------------------
$ws.Range("A1:B1").MergeCells = $true
$ws.Range("C1:D1").MergeCells = $true

# re-stamp the alignment so the merged cells pick back up the de-duped,
# centered style record instead of the blank one Merge() just created.
$ws.Range("A1:D1").HorizontalAlignment = -4108

# --- View state ------------------------------------------------------
$excel.ActiveWindow.Zoom = 93
[void]$ws.Range("D7").Select()
